# Organizational API Security Readiness workbook update
# - Rename "Sheet1" -> "Data" (defined name / autofilter reference updates automatically)
# - Insert a new header row on the "Legend" sheet and turn the A:B range into an
#   Excel Table ("Table1") with header labels "Column1"/"Column2"
# - Make the "Legend" sheet the active/selected tab, with A1:B7 selected
# - Remove the "tabSelected" flag from the Data sheet (handled automatically by activating Legend)

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsLegend = $wb.Worksheets.Item(2)

# Rename first sheet from "Sheet1" to "Data"
$wsData.Name = "Data"

# Insert a new row 1 on the Legend sheet to hold the table header
$wsLegend.Rows.Item(1).Insert()
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

# Give the new header columns a little extra width
$wsLegend.Columns.Item(1).ColumnWidth = 10.5
$wsLegend.Columns.Item(2).ColumnWidth = 10.5

# Make Legend the active sheet/tab
$wsLegend.Activate()

# Turn the A1:B7 range into a proper Excel Table
$lo = $wsLegend.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $wsLegend.Range("A1:B7"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Select the full table range as the active selection
[void]$wsLegend.Range("A1:B7").Select()
